$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 5: B5 changes from "en proceso" to "terminado"
$ws.Range("B5").Value = "terminado"

# New rows 6-8 (shared-string creation order matches the target file: ZZZ text,
# then "crear formulario...", then the row-5 comment, then "generar reporte...")
$ws.Range("A6").Value = "revisar formularios y permitir ingreso de codigo ZZZ"
$ws.Range("B6").Value = "en proceso"

$ws.Range("A7").Value = "crear formulario para cargar parametros del sistema"
$ws.Range("B7").Value = "no comenzado"

$ws.Range("C5").Value = "las consultas funcionan y mientras sigan asi no se ajustaran"

$ws.Range("A8").Value = "generar reporte de historico de movimientos de articulos"
$ws.Range("B8").Value = "en proceso"

# Column width adjustments (auto best-fit style widening of A and new col C)
$ws.Columns.Item(1).ColumnWidth = 71.9
$ws.Columns.Item(3).ColumnWidth = 53

# Update selection to match final state
$ws.Range("B11").Select()
